$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5929
$ws1.Range("F6").Value = 2980
$ws1.Range("F12").Value = 715
$ws1.Range("F13").Value = 250
$ws1.Range("F14").Value = 4340
$ws1.Range("F15").Value = 4340
$ws1.Range("F18").Value = 113
$ws1.Range("F19").Value = 23
$ws1.Range("F22").Value = 6636
$ws1.Range("F26").Value = 457
$ws1.Range("F27").Value = 1243
$ws1.Range("F28").Value = 6250
$ws1.Range("F32").Value = 5983
$ws1.Range("F35").Value = 95
$ws1.Range("F37").Value = 418
$ws1.Range("F38").Value = 4123
$ws1.Range("F40").Value = 188
$ws1.Range("F41").Value = 83
$ws1.Range("F49").Value = 2053

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 1
$ws2.Range("F13").Value = 19

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1415

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1415
$ws4.Range("F5").Value = 5929
$ws4.Range("F7").Value = 2980
$ws4.Range("F14").Value = 250
$ws4.Range("F15").Value = 4340
$ws4.Range("F16").Value = 4340
$ws4.Range("F19").Value = 113
$ws4.Range("F20").Value = 23
$ws4.Range("F22").Value = 6636
$ws4.Range("F25").Value = 457
$ws4.Range("F26").Value = 1243
$ws4.Range("F28").Value = 6250
$ws4.Range("F33").Value = 5983
$ws4.Range("F36").Value = 95
$ws4.Range("F38").Value = 418
$ws4.Range("F39").Value = 4123
$ws4.Range("F41").Value = 188
$ws4.Range("F42").Value = 83
$ws4.Range("F50").Value = 19

Write-Output "Done applying changes"
